$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 62.15419033333333
$ws.Range("H2").Value = 186.462571
$ws.Range("I2").Value = 0.5307382952913039
$ws.Range("J2").Value = 0.5798280707535227
$ws.Range("M2").Value = 90.43008666666667
$ws.Range("N2").Value = 271.29026
$ws.Range("O2").Value = 0.863466363695901
$ws.Range("P2").Value = 0.8656179140344247
$ws.Range("Q2").Value = 5620.608818539828
$ws.Range("R2").Value = 50585.47936685846
$ws.Range("S2").Value = 0.4582746659093435
$ws.Range("T2").Value = 0.5019095651042691
$ws.Range("G3").Value = 62.15419033333333
$ws.Range("H3").Value = 186.462571
$ws.Range("I3").Value = 0.5307382952913039
$ws.Range("J3").Value = 0.5798280707535227
$ws.Range("O3").Value = 0.000602137432244878
$ws.Range("P3").Value = 0.0006036378137891445
$ws.Range("Q3").Value = 3.919526114673777
$ws.Range("R3").Value = 35.275735032064
$ws.Range("S3").Value = 0.0003195773943207295
$ws.Range("T3").Value = 0.0003500061490032339
$ws.Range("G4").Value = 62.15419033333333
$ws.Range("H4").Value = 186.462571
$ws.Range("I4").Value = 0.5307382952913039
$ws.Range("J4").Value = 0.5798280707535227
$ws.Range("M4").Value = 9.467965
$ws.Range("N4").Value = 28.403895
$ws.Range("O4").Value = 0.0904043069236993
$ws.Range("P4").Value = 0.09062957269587499
$ws.Range("Q4").Value = 588.4736986793383
$ws.Range("R4").Value = 5296.263288114044
$ws.Range("S4").Value = 0.04798102774367599
$ws.Range("T4").Value = 0.05254957028946533
$ws.Range("G5").Value = 62.15419033333333
$ws.Range("H5").Value = 186.462571
$ws.Range("I5").Value = 0.5307382952913039
$ws.Range("J5").Value = 0.5798280707535227
$ws.Range("M5").Value = 0.7809334999999999
$ws.Range("N5").Value = 1.561867
$ws.Range("O5").Value = 0.007456697592460336
$ws.Range("P5").Value = 0.004983518592002547
$ws.Range("Q5").Value = 48.53828939667616
$ws.Range("R5").Value = 291.229736380057
$ws.Range("S5").Value = 0.003957554968725169
$ws.Range("T5").Value = 0.002889583970765149
$ws.Range("G6").Value = 62.15419033333333
$ws.Range("H6").Value = 186.462571
$ws.Range("I6").Value = 0.5307382952913039
$ws.Range("J6").Value = 0.5798280707535227
$ws.Range("M6").Value = 3.98709
$ws.Range("N6").Value = 11.96127
$ws.Range("O6").Value = 0.03807049435569441
$ws.Range("P6").Value = 0.03816535686390858
$ws.Range("Q6").Value = 247.81435073613
$ws.Range("R6").Value = 2230.32915662517
$ws.Range("S6").Value = 0.02020546927523846
$ws.Range("T6").Value = 0.02212934524001983
$ws.Range("I7").Value = 0.2152737834352902
$ws.Range("J7").Value = 0.2351851819258409
$ws.Range("M7").Value = 90.43008666666667
$ws.Range("N7").Value = 271.29026
$ws.Range("O7").Value = 0.863466363695901
$ws.Range("P7").Value = 0.8656179140344247
$ws.Range("Q7").Value = 2279.785981738354
$ws.Range("R7").Value = 20518.07383564518
$ws.Range("S7").Value = 0.1858816709819289
$ws.Range("T7").Value = 0.2035805065904531
$ws.Range("I8").Value = 0.2152737834352902
$ws.Range("J8").Value = 0.2351851819258409
$ws.Range("O8").Value = 0.000602137432244878
$ws.Range("P8").Value = 0.0006036378137891445
$ws.Range("S8").Value = 0.0001296244031873656
$ws.Range("T8").Value = 0.0001419666690533168
$ws.Range("I9").Value = 0.2152737834352902
$ws.Range("J9").Value = 0.2351851819258409
$ws.Range("M9").Value = 9.467965
$ws.Range("N9").Value = 28.403895
$ws.Range("O9").Value = 0.0904043069236993
$ws.Range("P9").Value = 0.09062957269587499
$ws.Range("Q9").Value = 238.691951741165
$ws.Range("R9").Value = 2148.227565670485
$ws.Range("S9").Value = 0.01946167719030995
$ws.Range("T9").Value = 0.02131473254234058
$ws.Range("I10").Value = 0.2152737834352902
$ws.Range("J10").Value = 0.2351851819258409
$ws.Range("M10").Value = 0.7809334999999999
$ws.Range("N10").Value = 1.561867
$ws.Range("O10").Value = 0.007456697592460336
$ws.Range("P10").Value = 0.004983518592002547
$ws.Range("Q10").Value = 19.6877091640135
$ws.Range("R10").Value = 118.126254984081
$ws.Range("S10").Value = 0.001605231502661756
$ws.Range("T10").Value = 0.00117204972669093
$ws.Range("I11").Value = 0.2152737834352902
$ws.Range("J11").Value = 0.2351851819258409
$ws.Range("M11").Value = 3.98709
$ws.Range("N11").Value = 11.96127
$ws.Range("O11").Value = 0.03807049435569441
$ws.Range("P11").Value = 0.03816535686390858
$ws.Range("Q11").Value = 100.51645669029
$ws.Range("R11").Value = 904.6481102126102
$ws.Range("S11").Value = 0.008195579357202196
$ws.Range("T11").Value = 0.008975926397302981
$ws.Range("G12").Value = 29.744252
$ws.Range("H12").Value = 59.488504
$ws.Range("I12").Value = 0.2539879212734059
$ws.Range("J12").Value = 0.1849867473206364
$ws.Range("M12").Value = 90.43008666666667
$ws.Range("N12").Value = 271.29026
$ws.Range("O12").Value = 0.863466363695901
$ws.Range("P12").Value = 0.8656179140344247
$ws.Range("Q12").Value = 2689.775286195174
$ws.Range("R12").Value = 16138.65171717104
$ws.Range("S12").Value = 0.2193100268046286
$ws.Range("T12").Value = 0.1601278423397025
$ws.Range("G13").Value = 29.744252
$ws.Range("H13").Value = 59.488504
$ws.Range("I13").Value = 0.2539879212734059
$ws.Range("J13").Value = 0.1849867473206364
$ws.Range("O13").Value = 0.000602137432244878
$ws.Range("P13").Value = 0.0006036378137891445
$ws.Range("Q13").Value = 1.875712190122667
$ws.Range("R13").Value = 11.254273140736
$ws.Range("S13").Value = 0.0001529356347367829
$ws.Range("T13").Value = 0.0001116649957325938
$ws.Range("G14").Value = 29.744252
$ws.Range("H14").Value = 59.488504
$ws.Range("I14").Value = 0.2539879212734059
$ws.Range("J14").Value = 0.1849867473206364
$ws.Range("M14").Value = 9.467965
$ws.Range("N14").Value = 28.403895
$ws.Range("O14").Value = 0.0904043069236993
$ws.Range("P14").Value = 0.09062957269587499
$ws.Range("Q14").Value = 281.61753688718
$ws.Range("R14").Value = 1689.70522132308
$ws.Range("S14").Value = 0.02296160198971336
$ws.Range("T14").Value = 0.01676526986406907
$ws.Range("G15").Value = 29.744252
$ws.Range("H15").Value = 59.488504
$ws.Range("I15").Value = 0.2539879212734059
$ws.Range("J15").Value = 0.1849867473206364
$ws.Range("M15").Value = 0.7809334999999999
$ws.Range("N15").Value = 1.561867
$ws.Range("O15").Value = 0.007456697592460336
$ws.Range("P15").Value = 0.004983518592002547
$ws.Range("Q15").Value = 23.228282819242
$ws.Range("R15").Value = 92.91313127696799
$ws.Range("S15").Value = 0.001893911121073411
$ws.Range("T15").Value = 0.0009218848945464688
$ws.Range("G16").Value = 29.744252
$ws.Range("H16").Value = 59.488504
$ws.Range("I16").Value = 0.2539879212734059
$ws.Range("J16").Value = 0.1849867473206364
$ws.Range("M16").Value = 3.98709
$ws.Range("N16").Value = 11.96127
$ws.Range("O16").Value = 0.03807049435569441
$ws.Range("P16").Value = 0.03816535686390858
$ws.Range("Q16").Value = 118.59300970668
$ws.Range("R16").Value = 711.5580582400801
$ws.Range("S16").Value = 0.009669445723253758
$ws.Range("T16").Value = 0.02212934524001983
